# Update cryptocurrency price/volume figures (refreshed data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.118.36"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.885.02"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.05"
$ws.Range("E5").Value = "  -3.03%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4701"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4021"
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.35"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08000"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9921"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.37"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "1.904.70"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.867"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.016"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.69"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06607"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001021"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "29.139.27"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.479"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.61"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.180"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "2.108.58"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.98"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.017"
$ws.Range("E29").Value = "  +6.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.067"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.31"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.032"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09421"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.539"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.378"
$ws.Range("E35").Value = "  -3.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.338"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06060"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02223"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.170"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.008"
$ws.Range("E40").Value = "  -5.37%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1823"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.448"
$ws.Range("E43").Value = "  +5.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.988"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07685"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.17"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5461"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.79"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.38"
$ws.Range("E51").Value = "  +0.70%  "
